$p = $ppt.ActivePresentation
try {
  $cp = $p.CustomXMLParts.Add("<a:theme xmlns:a='http://schemas.openxmlformats.org/drawingml/2006/main'/>")
  Write-Host "added customxml: $cp"
} catch {
  Write-Host "ERROR: $($_.Exception.Message)"
}
